$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7832693457603455
$ws.Range("B1").Value = 3.380542516708374
$ws.Range("C1").Value = 3.547896862030029
$ws.Range("D1").Value = 2.932676792144775
$ws.Range("E1").Value = 1.781857132911682
